# Insert a new weekly record before the current row 230, shifting all
# subsequent rows (230..280) down by one (to 231..281). The new row 230
# duplicates the data that was previously in row 230, except for a new
# date (Fecha) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 230; Excel shifts rows 230:280 down
# to 231:281 (formats/styles carry down with the existing cells).
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the record (same as the
# original row 230's data, but with the new Fecha).
$ws.Range("A230").Value2 = 1
$ws.Range("B230").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C230").Value2 = "Arica y Parinacota"
$ws.Range("D230").Value2 = 44798
$ws.Range("E230").Value2 = 15
$ws.Range("F230").Value2 = "Fruta"
$ws.Range("G230").Value2 = 100108
$ws.Range("H230").Value2 = "Tropicales y subtropicales"
$ws.Range("I230").Value2 = 100108006
$ws.Range("J230").Value2 = "Plátano"
$ws.Range("K230").Value2 = "Sin especificar"
$ws.Range("L230").Value2 = "Pintón"
$ws.Range("M230").Value2 = 120
$ws.Range("N230").Value2 = 19000
$ws.Range("O230").Value2 = 20000
$ws.Range("P230").Value2 = 19500
$ws.Range("Q230").Value2 = "$/caja 20 kilos"
$ws.Range("R230").Value2 = "Ecuador"
$ws.Range("S230").Value2 = 975
$ws.Range("T230").Value2 = 20
